# Weekly data update: insert a new record as row 41 (pushing every
# subsequent "Haba" price row down by one), matching the site's
# weekly "logica_diaria" refresh pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41..62 down to 42..63, opening up a blank row 41.
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with this week's record.
$ws.Cells.Item(41, 1).Value = 2
$ws.Cells.Item(41, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(41, 3).Value = "Coquimbo"
$ws.Cells.Item(41, 4).Value = 44777
$ws.Cells.Item(41, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41, 5).Value = 4
$ws.Cells.Item(41, 6).Value = 100112026
$ws.Cells.Item(41, 7).Value = "Haba"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 600
$ws.Cells.Item(41, 11).Value = 9000
$ws.Cells.Item(41, 12).Value = 10000
$ws.Cells.Item(41, 13).Value = 9500
$ws.Cells.Item(41, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(41, 16).Value = 380
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
